# Actualización automática 2025-07-03 15:05:08
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("H19").Value = 426.6
$ws1.Range("M35").Value = 6579.28
$ws1.Range("H55").Value = "1 de 53"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F19").Value = 426.6
$ws2.Range("F35").Value = 6732.78
$ws2.Range("F55").Value = 10052.26

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Columns.Item(6).ColumnWidth = 23.17

$ws3.Range("D7").Value = 426.6
$ws3.Range("E7").Value = 1973.4
$ws3.Range("F7").Value = 0.17775

$ws3.Range("D16").Value = 9472.16
$ws3.Range("E16").Value = 42354.3
$ws3.Range("F16").Value = 0.1827668723659691

$ws3.Range("D19").Value = 10052.26
$ws3.Range("E19").Value = 103654.1906451792
$ws3.Range("F19").Value = 0.08840536260663054
